$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 858.4167
$ws.Range("I98").Value = 858.4167
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 858.4167
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 639.5833
$ws.Range("N98").ClearContents()

$ws.Range("H112").Value = 967.7778
$ws.Range("J112").Value = 1168.3334
$ws.Range("L112").Value = 3505.0002
$ws.Range("N112").Value = -5721.0002

$ws.Range("H122").Value = 858.4167
$ws.Range("I122").Value = 858.4167
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2575.2501
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -125.2501000000002
$ws.Range("N122").ClearContents()

$ws.Range("H129").Value = 885.5
$ws.Range("J129").Value = 1075.9459
$ws.Range("L129").Value = 3227.8377
$ws.Range("N129").Value = -13227.8377

$ws.Range("H132").Value = 3224.648
$ws.Range("I132").Value = 3153.3333
$ws.Range("K132").Value = 9459.999899999999
$ws.Range("M132").Value = -6929.999899999999

$ws.Range("H138").Value = 2905.0266
$ws.Range("I138").Value = 1237.0741
$ws.Range("J138").Value = 3843.25
$ws.Range("K138").Value = 3711.2223
$ws.Range("L138").Value = 11529.75
$ws.Range("M138").Value = 1428.7777
$ws.Range("N138").Value = -21809.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1210.0385
$ws.Range("I2").Value = 1076.8948
$ws.Range("J2").Value = 1571.4286
$ws.Range("K2").Value = 1076.8948
$ws.Range("L2").Value = 1571.4286
$ws.Range("M2").Value = -963.8948
$ws.Range("N2").Value = -1797.4286

$ws.Range("H61").Value = 1628.279
$ws.Range("I61").Value = 1379.8334
$ws.Range("J61").Value = 2906
$ws.Range("K61").Value = 1379.8334
$ws.Range("L61").Value = 2906
$ws.Range("M61").Value = -1167.8334
$ws.Range("N61").Value = -3330

$ws.Range("H74").Value = 224247.47
$ws.Range("I74").Value = 2276.7058
$ws.Range("J74").Value = 910338.9399999999
$ws.Range("K74").Value = 2276.7058
$ws.Range("L74").Value = 910338.9399999999
$ws.Range("M74").Value = -1402.7058
$ws.Range("N74").Value = -912086.9399999999

$ws.Range("H77").Value = 224247.47
$ws.Range("I77").Value = 2276.7058
$ws.Range("J77").Value = 910338.9399999999
$ws.Range("K77").Value = 11383.529
$ws.Range("L77").Value = 4551694.699999999
$ws.Range("M77").Value = -7015.529
$ws.Range("N77").Value = -4560430.699999999

$ws.Range("H110").Value = 2725.125
$ws.Range("I110").Value = 2841.8333
$ws.Range("J110").Value = 2375
$ws.Range("K110").Value = 2841.8333
$ws.Range("L110").Value = 2375
$ws.Range("M110").Value = -796.8332999999998
$ws.Range("N110").Value = -6465

$ws.Range("H116").Value = 1210.0385
$ws.Range("I116").Value = 1076.8948
$ws.Range("J116").Value = 1571.4286
$ws.Range("K116").Value = 1076.8948
$ws.Range("L116").Value = 1571.4286
$ws.Range("M116").Value = 1217.1052
$ws.Range("N116").Value = -6159.4286

$ws.Range("H132").Value = 18552.17
$ws.Range("I132").Value = 23966.158
$ws.Range("J132").Value = 2671.1333
$ws.Range("K132").Value = 71898.474
$ws.Range("L132").Value = 8013.3999
$ws.Range("M132").Value = -69368.474
$ws.Range("N132").Value = -13073.3999

$ws.Range("H136").Value = 1628.279
$ws.Range("I136").Value = 1379.8334
$ws.Range("J136").Value = 2906
$ws.Range("K136").Value = 4139.5002
$ws.Range("L136").Value = 8718
$ws.Range("M136").Value = -1589.5002
$ws.Range("N136").Value = -13818

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1210.0385
$ws.Range("I3").Value = 1076.8948
$ws.Range("J3").Value = 1571.4286
$ws.Range("K3").Value = 1076.8948
$ws.Range("L3").Value = 1571.4286
$ws.Range("M3").Value = -962.8948
$ws.Range("N3").Value = -1799.4286

$ws.Range("H86").Value = 4652884
$ws.Range("I86").Value = 7144585
$ws.Range("J86").Value = 1709.2667
$ws.Range("K86").Value = 7144585
$ws.Range("L86").Value = 1709.2667
$ws.Range("M86").Value = -7143462
$ws.Range("N86").Value = -3955.2667

$ws.Range("H89").Value = 4652884
$ws.Range("I89").Value = 7144585
$ws.Range("J89").Value = 1709.2667
$ws.Range("K89").Value = 35722925
$ws.Range("L89").Value = 8546.333499999999
$ws.Range("M89").Value = -35717309
$ws.Range("N89").Value = -19778.3335

$ws.Range("H107").Value = 5381.154
$ws.Range("I107").Value = 6734.316
$ws.Range("J107").Value = 1708.2858
$ws.Range("K107").Value = 6734.316
$ws.Range("L107").Value = 1708.2858
$ws.Range("M107").Value = -4814.316
$ws.Range("N107").Value = -5548.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1155.4615
$ws.Range("I16").Value = 974.63635
$ws.Range("J16").Value = 2150
$ws.Range("K16").Value = 974.63635
$ws.Range("L16").Value = 2150
$ws.Range("M16").Value = -687.63635
$ws.Range("N16").Value = -2724

$ws.Range("H31").Value = 1936.791
$ws.Range("I31").Value = 1165.3143
$ws.Range("J31").Value = 2780.5938
$ws.Range("K31").Value = 1165.3143
$ws.Range("L31").Value = 2780.5938
$ws.Range("M31").Value = -870.3143
$ws.Range("N31").Value = -3370.5938

$ws.Range("H34").Value = 1936.791
$ws.Range("I34").Value = 1165.3143
$ws.Range("J34").Value = 2780.5938
$ws.Range("K34").Value = 1165.3143
$ws.Range("L34").Value = 2780.5938
$ws.Range("M34").Value = -963.3143
$ws.Range("N34").Value = -3184.5938

$ws.Range("H113").Value = 1155.4615
$ws.Range("I113").Value = 974.63635
$ws.Range("J113").Value = 2150
$ws.Range("K113").Value = 974.63635
$ws.Range("L113").Value = 2150
$ws.Range("M113").Value = 1195.36365
$ws.Range("N113").Value = -6490

$ws.Range("H122").Value = 1110.826
$ws.Range("I122").Value = 857.61536
$ws.Range("J122").Value = 1440
$ws.Range("K122").Value = 2572.84608
$ws.Range("L122").Value = 4320
$ws.Range("M122").Value = -122.8460800000003
$ws.Range("N122").Value = -9220

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 716.60376
$ws.Range("I5").Value = 411.39285
$ws.Range("J5").Value = 1058.44
$ws.Range("K5").Value = 1234.17855
$ws.Range("L5").Value = 3175.32
$ws.Range("M5").Value = -1122.17855
$ws.Range("N5").Value = -3399.32

$ws.Range("H122").Value = 775.5
$ws.Range("I122").Value = 447.0476
$ws.Range("J122").Value = 1103.9524
$ws.Range("K122").Value = 4023.4284
$ws.Range("L122").Value = 9935.571599999999
$ws.Range("M122").Value = -1573.4284
$ws.Range("N122").Value = -14835.5716

$ws.Range("H131").Value = 973
$ws.Range("J131").Value = 973
$ws.Range("L131").Value = 2919
$ws.Range("N131").Value = -12999

$ws.Range("H133").Value = 981.0625
$ws.Range("I133").Value = 449.5
$ws.Range("K133").Value = 1348.5
$ws.Range("M133").Value = 3711.5

$ws.Range("H134").Value = 1084.3334
$ws.Range("I134").Value = 814.75
$ws.Range("J134").Value = 1300
$ws.Range("K134").Value = 2444.25
$ws.Range("L134").Value = 3900
$ws.Range("M134").Value = 2625.75
$ws.Range("N134").Value = -14040

$ws.Range("H135").Value = 716.60376
$ws.Range("I135").Value = 411.39285
$ws.Range("J135").Value = 1058.44
$ws.Range("K135").Value = 3702.53565
$ws.Range("L135").Value = 9525.960000000001
$ws.Range("M135").Value = -1167.53565
$ws.Range("N135").Value = -14595.96

$ws.Range("H136").Value = 1166.3572

$ws.Range("H137").Value = 2925.611
$ws.Range("I137").Value = 951.26666
$ws.Range("J137").Value = 4335.857
$ws.Range("K137").Value = 2853.79998
$ws.Range("L137").Value = 13007.571
$ws.Range("M137").Value = 2246.20002
$ws.Range("N137").Value = -23207.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 63805452
$ws.Range("I70").Value = 138237700
$ws.Range("K70").Value = 138237700
$ws.Range("M70").Value = -138237430

$ws.Range("H73").Value = 63805452
$ws.Range("I73").Value = 138237700
$ws.Range("K73").Value = 138237700
$ws.Range("M73").Value = -138236764

$ws.Range("H132").Value = 1933.898
$ws.Range("I132").Value = 1320.9678
$ws.Range("J132").Value = 2989.5
$ws.Range("K132").Value = 3962.9034
$ws.Range("L132").Value = 8968.5
$ws.Range("M132").Value = -1432.9034
$ws.Range("N132").Value = -14028.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1926.6
$ws.Range("I61").Value = 1607.2307
$ws.Range("J61").Value = 4002.5
$ws.Range("K61").Value = 1607.2307
$ws.Range("L61").Value = 4002.5
$ws.Range("M61").Value = -1405.2307
$ws.Range("N61").Value = -4406.5

$ws.Range("H113").Value = 1926.6
$ws.Range("I113").Value = 1607.2307
$ws.Range("J113").Value = 4002.5
$ws.Range("K113").Value = 1607.2307
$ws.Range("L113").Value = 4002.5
$ws.Range("M113").Value = 562.7692999999999
$ws.Range("N113").Value = -8342.5

$ws.Range("H122").Value = 4933.1665
$ws.Range("I122").Value = 4995.273
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 14985.819
$ws.Range("L122").Value = 12750
$ws.Range("M122").Value = -12535.819
$ws.Range("N122").Value = -17650

$ws.Range("H132").Value = 9264960
$ws.Range("I132").Value = 11910936
$ws.Range("J132").Value = 4044.8333
$ws.Range("K132").Value = 35732808
$ws.Range("L132").Value = 12134.4999
$ws.Range("M132").Value = -35730278
$ws.Range("N132").Value = -17194.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1134.0358
$ws.Range("I107").Value = 1026.4667
$ws.Range("J107").Value = 1258.1538
$ws.Range("K107").Value = 3079.4001
$ws.Range("L107").Value = 3774.4614
$ws.Range("M107").Value = -1159.4001
$ws.Range("N107").Value = -7614.4614

$ws.Range("H122").Value = 948.64514
$ws.Range("I122").Value = 816
$ws.Range("J122").Value = 1272.8889
$ws.Range("K122").Value = 2448
$ws.Range("L122").Value = 3818.6667
$ws.Range("M122").Value = 2
$ws.Range("N122").Value = -8718.6667

$ws.Range("H132").Value = 1565.9762
$ws.Range("I132").Value = 842.8823
$ws.Range("J132").Value = 4639.125
$ws.Range("K132").Value = 2528.6469
$ws.Range("L132").Value = 13917.375
$ws.Range("M132").Value = 1.353100000000268
$ws.Range("N132").Value = -18977.375

$ws.Range("H136").Value = 4174.575
$ws.Range("I136").Value = 1213.7142
$ws.Range("J136").Value = 11083.25
$ws.Range("K136").Value = 3641.1426
$ws.Range("L136").Value = 33249.75
$ws.Range("M136").Value = -1091.1426
$ws.Range("N136").Value = -38349.75

$ws.Range("H138").Value = 51050
$ws.Range("I138").Value = 12000
$ws.Range("J138").Value = 58860
$ws.Range("K138").Value = 12000
$ws.Range("L138").Value = 58860
$ws.Range("N138").Value = -69140
$ws.Range("M138").Value = -6860
